$wb = $excel.ActiveWorkbook

# ---- Sheet: 展览 ----
$ws = $wb.Worksheets.Item("展览")
$ws.Range("F4").Value  = 99
$ws.Range("F5").Value  = 63
$ws.Range("F6").Value  = 700
$ws.Range("F8").Value  = 198
$ws.Range("F9").Value  = 14
$ws.Range("F12").Value = 561
$ws.Range("F13").Value = 483
$ws.Range("F16").Value = 135
$ws.Range("F17").Value = 777
$ws.Range("F18").Value = 2578
$ws.Range("F23").Value = 187
$ws.Range("F24").Value = 13
$ws.Range("F25").Value = 132
$ws.Range("F26").Value = 576
$ws.Range("F27").Value = 943
$ws.Range("F29").Value = 156
$ws.Range("F30").Value = 68
$ws.Range("F33").Value = 256

# ---- Sheet: 演出 ----
$ws = $wb.Worksheets.Item("演出")
$ws.Range("F4").Value  = 333
$ws.Range("F9").Value  = 293
$ws.Range("F13").Value = 526
$ws.Range("F16").Value = 953
$ws.Range("F23").Value = 262
$ws.Range("F24").Value = 236

# ---- Sheet: 本地生活 ----
$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F3").Value  = 74
$ws.Range("F5").Value  = 2351
$ws.Range("F6").Value  = 950
$ws.Range("G8").Value  = "不可售"

# ---- Sheet: 全部类型 ----
$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F3").Value  = 74
$ws.Range("F5").Value  = 2351
$ws.Range("F9").Value  = 950
$ws.Range("F13").Value = 99
$ws.Range("F14").Value = 63
$ws.Range("F15").Value = 700
$ws.Range("F18").Value = 198
$ws.Range("F19").Value = 14
$ws.Range("F21").Value = 561
$ws.Range("F22").Value = 483
$ws.Range("F25").Value = 135
$ws.Range("F26").Value = 777
$ws.Range("F27").Value = 2578
$ws.Range("F31").Value = 293
$ws.Range("F32").Value = 187
$ws.Range("F33").Value = 132
$ws.Range("F34").Value = 576
$ws.Range("F35").Value = 943
$ws.Range("F36").Value = 526
$ws.Range("F39").Value = 156
$ws.Range("F44").Value = 262
$ws.Range("F45").Value = 236
$ws.Range("F49").Value = 256

$wb.Save()
